$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Trip Costs sheet: fill in the "Total" column (D) formulas and the
# "Amount per student" cell.
# ---------------------------------------------------------------------------
$wsTrip = $wb.Worksheets.Item("Trip Costs")

$wsTrip.Range("D6").Formula  = "=B6*`$B`$3"
$wsTrip.Range("D7").Formula  = "=B7*C7"
$wsTrip.Range("D8").Formula  = "=B8*`$B`$3"
$wsTrip.Range("D9").Formula  = "=B9*`$B`$3"
$wsTrip.Range("D10").Formula = "=B10*`$B`$3"
$wsTrip.Range("D11").Formula = "=B11*C11"
$wsTrip.Range("D12").Formula = "=SUM(D6:D11)"
$wsTrip.Range("D14").Formula = "=D12/B3"

# ---------------------------------------------------------------------------
# Payments sheet: fill in "Total Paid" (G) and "Outstanding" (H) columns for
# every student row, plus the TOTALS row.
# ---------------------------------------------------------------------------
$wsPay = $wb.Worksheets.Item("Payments")

for ($r = 4; $r -le 15; $r++) {
    $wsPay.Range("G$r").Formula = "=SUM(D$r`:F$r)"
    $wsPay.Range("H$r").Formula = "=C$r-G$r"
}

$wsPay.Range("C16").Formula = "=SUM(C4:C15)"
$wsPay.Range("G16").Formula = "=SUM(G4:G15)"
$wsPay.Range("H16").Formula = "=SUM(H4:H15)"

# ---------------------------------------------------------------------------
# View state: update selections on each sheet and make "Payments" the active
# (visible) tab, matching where the author left off working.
# ---------------------------------------------------------------------------
$wsTitle = $wb.Worksheets.Item("Title page")
$wsTitle.Activate() | Out-Null
$wsTitle.Range("A12:G12").Select() | Out-Null

$wsTrip.Activate() | Out-Null
$wsTrip.Range("D15").Select() | Out-Null

$wsPay.Activate() | Out-Null
$wsPay.Range("G15").Select() | Out-Null
